$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A5").Value = 42606.881793981483
$ws.Range("A5").NumberFormat = "m/d/yy h:mm"

$ws.Range("B5").Value = 26
$ws.Range("C5").Value = 64
$ws.Range("D5").Value = 34
$ws.Range("E5").Value = 60
$ws.Range("F5").Value = 40
$ws.Range("G5").Value = 10646
$ws.Range("H5").Value = 8010
$ws.Range("I5").Value = 486
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 54
$ws.Range("L5").Value = 3
$ws.Range("M5").Value = 2
$ws.Range("N5").Value = "Named"
